# Insert a new arrival row (18:47 14X44_ABASTO 0 🚌) as row 3 in both the
# "TODOS" and "COMBINADAS" sheets, shifting the existing rows 3..43 down by one.
$wb = $excel.ActiveWorkbook

$sheetNames = @("TODOS", "COMBINADAS")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Push existing data down by inserting a new row before row 3.
    $ws.Rows.Item(3).Insert()

    # Populate the newly inserted row.
    $ws.Cells.Item(3, 1).Value = "18:47"
    $ws.Cells.Item(3, 2).Value = "14X44_ABASTO"
    $ws.Cells.Item(3, 3).Value = 0
    $ws.Cells.Item(3, 4).Value = "🚌"
}
